$wb = $excel.ActiveWorkbook

# ALC (sheet1)
$ws = $wb.Worksheets.Item("ALC")
# row 9
$ws.Cells.Item(9, 8).Value = 92.25
$ws.Cells.Item(9, 9).Value = 73.166664
$ws.Cells.Item(9, 11).Value = 73.166664
$ws.Cells.Item(9, 13).Value = 95.833336
# row 28
$ws.Cells.Item(28, 8).Value = 166667400
$ws.Cells.Item(28, 9).Value = 166667400
$ws.Cells.Item(28, 11).Value = 166667400
$ws.Cells.Item(28, 13).Value = -166666915
# row 64
$ws.Cells.Item(64, 8).Value = 3800
$ws.Cells.Item(64, 9).Value = 5200
$ws.Cells.Item(64, 10).Value = 2400
$ws.Cells.Item(64, 11).Value = 5200
$ws.Cells.Item(64, 12).Value = 2400
$ws.Cells.Item(64, 13).Value = -4952
$ws.Cells.Item(64, 14).Value = -2896
# row 67
$ws.Cells.Item(67, 8).Value = 3800
$ws.Cells.Item(67, 9).Value = 5200
$ws.Cells.Item(67, 10).Value = 2400
$ws.Cells.Item(67, 11).Value = 5200
$ws.Cells.Item(67, 12).Value = 2400
$ws.Cells.Item(67, 13).Value = -4342
$ws.Cells.Item(67, 14).Value = -4116
# row 98
$ws.Cells.Item(98, 8).Value = 1815.9166
$ws.Cells.Item(98, 9).Value = 1859.3
$ws.Cells.Item(98, 11).Value = 1859.3
$ws.Cells.Item(98, 13).Value = -361.3
# row 122
$ws.Cells.Item(122, 8).Value = 1815.9166
$ws.Cells.Item(122, 9).Value = 1859.3
$ws.Cells.Item(122, 11).Value = 5577.9
$ws.Cells.Item(122, 13).Value = -3127.9
# row 138
$ws.Cells.Item(138, 8).Value = 3417.3547
$ws.Cells.Item(138, 9).Value = 1852.7142
$ws.Cells.Item(138, 10).Value = 3873.7083
$ws.Cells.Item(138, 11).Value = 5558.142599999999
$ws.Cells.Item(138, 12).Value = 11621.1249
$ws.Cells.Item(138, 13).Value = -418.1425999999992
$ws.Cells.Item(138, 14).Value = -21901.1249

# ARM (sheet2)
$ws = $wb.Worksheets.Item("ARM")
# row 10
$ws.Cells.Item(10, 8).Value = 2775
$ws.Cells.Item(10, 10).Value = 5500
$ws.Cells.Item(10, 12).Value = 5500
$ws.Cells.Item(10, 14).Value = -5840
# row 138
$ws.Cells.Item(138, 8).Value = 0
$ws.Cells.Item(138, 10).Value = 0
$ws.Cells.Item(138, 12).Value = 0
$ws.Cells.Item(138, 14).ClearContents()

# BSM (sheet3)
$ws = $wb.Worksheets.Item("BSM")
# row 19
$ws.Cells.Item(19, 8).Value = 11280.8
$ws.Cells.Item(19, 10).Value = 16999.334
$ws.Cells.Item(19, 12).Value = 16999.334
$ws.Cells.Item(19, 14).Value = -17345.334
# row 22
$ws.Cells.Item(22, 8).Value = 153.83333
$ws.Cells.Item(22, 9).Value = 60.75
$ws.Cells.Item(22, 11).Value = 60.75
$ws.Cells.Item(22, 13).Value = 112.25

# CRP (sheet4)
$ws = $wb.Worksheets.Item("CRP")
# row 25
$ws.Cells.Item(25, 8).Value = 1447.9166
# row 31
$ws.Cells.Item(31, 8).Value = 2281.2258
$ws.Cells.Item(31, 9).Value = 1309.72
$ws.Cells.Item(31, 11).Value = 1309.72
$ws.Cells.Item(31, 13).Value = -1014.72
# row 34
$ws.Cells.Item(34, 8).Value = 2281.2258
$ws.Cells.Item(34, 9).Value = 1309.72
$ws.Cells.Item(34, 11).Value = 1309.72
$ws.Cells.Item(34, 13).Value = -1107.72
# row 58
$ws.Cells.Item(58, 8).Value = 2815.8333
$ws.Cells.Item(58, 9).Value = 1805
$ws.Cells.Item(58, 11).Value = 1805
$ws.Cells.Item(58, 13).Value = -1602
# row 62
$ws.Cells.Item(62, 8).Value = 103514.5
$ws.Cells.Item(62, 9).Value = 4686.3335
$ws.Cells.Item(62, 11).Value = 4686.3335
$ws.Cells.Item(62, 13).Value = -4062.3335
# row 65
$ws.Cells.Item(65, 8).Value = 103514.5
$ws.Cells.Item(65, 9).Value = 4686.3335
$ws.Cells.Item(65, 11).Value = 23431.6675
$ws.Cells.Item(65, 13).Value = -20311.6675
# row 136
$ws.Cells.Item(136, 8).Value = 2815.8333
$ws.Cells.Item(136, 9).Value = 1805
$ws.Cells.Item(136, 11).Value = 5415
$ws.Cells.Item(136, 13).Value = -2865
# row 141
$ws.Cells.Item(141, 8).Value = 0
$ws.Cells.Item(141, 10).Value = 0
$ws.Cells.Item(141, 12).Value = 0
$ws.Cells.Item(141, 14).ClearContents()

# CUL (sheet5)
$ws = $wb.Worksheets.Item("CUL")
# row 68
$ws.Cells.Item(68, 8).Value = 510.5
$ws.Cells.Item(68, 10).Value = 749.5
$ws.Cells.Item(68, 12).Value = 2248.5
$ws.Cells.Item(68, 14).Value = -3870.5
# row 71
$ws.Cells.Item(71, 8).Value = 510.5
$ws.Cells.Item(71, 10).Value = 749.5
$ws.Cells.Item(71, 12).Value = 6745.5
$ws.Cells.Item(71, 14).Value = -14857.5

# GSM (sheet6)
$ws = $wb.Worksheets.Item("GSM")
# row 47
$ws.Cells.Item(47, 8).Value = 22479
$ws.Cells.Item(47, 10).Value = 20000.5
$ws.Cells.Item(47, 12).Value = 20000.5
$ws.Cells.Item(47, 14).Value = -21136.5
# row 55
$ws.Cells.Item(55, 8).Value = 3750
$ws.Cells.Item(55, 10).Value = 4000
$ws.Cells.Item(55, 12).Value = 4000
$ws.Cells.Item(55, 14).Value = -4654
# row 70
$ws.Cells.Item(70, 8).Value = 15495.167
$ws.Cells.Item(70, 9).Value = 5994
$ws.Cells.Item(70, 10).Value = 24996.334
$ws.Cells.Item(70, 11).Value = 5994
$ws.Cells.Item(70, 12).Value = 24996.334
$ws.Cells.Item(70, 13).Value = -5724
$ws.Cells.Item(70, 14).Value = -25536.334
# row 73
$ws.Cells.Item(73, 8).Value = 15495.167
$ws.Cells.Item(73, 9).Value = 5994
$ws.Cells.Item(73, 10).Value = 24996.334
$ws.Cells.Item(73, 11).Value = 5994
$ws.Cells.Item(73, 12).Value = 24996.334
$ws.Cells.Item(73, 13).Value = -5058
$ws.Cells.Item(73, 14).Value = -26868.334
# row 126
$ws.Cells.Item(126, 8).Value = 6315.8335
$ws.Cells.Item(126, 9).Value = 6322.5
$ws.Cells.Item(126, 11).Value = 18967.5
$ws.Cells.Item(126, 13).Value = -16497.5

# LTW (sheet7)
$ws = $wb.Worksheets.Item("LTW")
# row 21
$ws.Cells.Item(21, 8).Value = 6500
$ws.Cells.Item(21, 10).Value = 7142.857
$ws.Cells.Item(21, 12).Value = 7142.857
$ws.Cells.Item(21, 14).Value = -7490.857
# row 22
$ws.Cells.Item(22, 8).Value = 6753.577
$ws.Cells.Item(22, 9).Value = 2518.6155
$ws.Cells.Item(22, 10).Value = 10988.538
$ws.Cells.Item(22, 11).Value = 2518.6155
$ws.Cells.Item(22, 12).Value = 10988.538
$ws.Cells.Item(22, 13).Value = -2223.6155
$ws.Cells.Item(22, 14).Value = -11578.538
# row 27
$ws.Cells.Item(27, 8).Value = 6753.577
$ws.Cells.Item(27, 9).Value = 2518.6155
$ws.Cells.Item(27, 10).Value = 10988.538
$ws.Cells.Item(27, 11).Value = 2518.6155
$ws.Cells.Item(27, 12).Value = 10988.538
$ws.Cells.Item(27, 13).Value = -2411.6155
$ws.Cells.Item(27, 14).Value = -11202.538
# row 46
$ws.Cells.Item(46, 8).Value = 4811.25
$ws.Cells.Item(46, 9).Value = 3000
$ws.Cells.Item(46, 10).Value = 5415
$ws.Cells.Item(46, 11).Value = 3000
$ws.Cells.Item(46, 12).Value = 5415
$ws.Cells.Item(46, 13).Value = -2812
$ws.Cells.Item(46, 14).Value = -5791
# row 61
$ws.Cells.Item(61, 8).Value = 18521618
$ws.Cells.Item(61, 9).Value = 22225342
$ws.Cells.Item(61, 11).Value = 22225342
$ws.Cells.Item(61, 13).Value = -22225140
# row 68
$ws.Cells.Item(68, 8).Value = 3749.25
$ws.Cells.Item(68, 9).Value = 3332.6667
$ws.Cells.Item(68, 10).Value = 4999
$ws.Cells.Item(68, 11).Value = 3332.6667
$ws.Cells.Item(68, 12).Value = 4999
$ws.Cells.Item(68, 13).Value = -2583.6667
$ws.Cells.Item(68, 14).Value = -6497
# row 71
$ws.Cells.Item(71, 8).Value = 3749.25
$ws.Cells.Item(71, 9).Value = 3332.6667
$ws.Cells.Item(71, 10).Value = 4999
$ws.Cells.Item(71, 11).Value = 16663.3335
$ws.Cells.Item(71, 12).Value = 24995
$ws.Cells.Item(71, 13).Value = -12919.3335
$ws.Cells.Item(71, 14).Value = -32483
# row 100
$ws.Cells.Item(100, 8).Value = 6706
$ws.Cells.Item(100, 9).Value = 6857.5
$ws.Cells.Item(100, 11).Value = 6857.5
$ws.Cells.Item(100, 13).Value = -6316.5
# row 113
$ws.Cells.Item(113, 8).Value = 18521618
$ws.Cells.Item(113, 9).Value = 22225342
$ws.Cells.Item(113, 11).Value = 22225342
$ws.Cells.Item(113, 13).Value = -22223172
# row 132
$ws.Cells.Item(132, 8).Value = 65031.625
$ws.Cells.Item(132, 9).Value = 85567.25
$ws.Cells.Item(132, 10).Value = 3424.75
$ws.Cells.Item(132, 11).Value = 256701.75
$ws.Cells.Item(132, 12).Value = 10274.25
$ws.Cells.Item(132, 13).Value = -254171.75
$ws.Cells.Item(132, 14).Value = -15334.25
# row 136
$ws.Cells.Item(136, 8).Value = 6887.8335
$ws.Cells.Item(136, 9).Value = 6443.3335
$ws.Cells.Item(136, 11).Value = 19330.0005
$ws.Cells.Item(136, 13).Value = -16780.0005

# WVR (sheet8)
$ws = $wb.Worksheets.Item("WVR")
# row 14
$ws.Cells.Item(14, 8).Value = 18998.166
$ws.Cells.Item(14, 9).Value = 21000
$ws.Cells.Item(14, 10).Value = 18597.8
$ws.Cells.Item(14, 11).Value = 21000
$ws.Cells.Item(14, 12).Value = 18597.8
$ws.Cells.Item(14, 13).Value = -20832
$ws.Cells.Item(14, 14).Value = -18933.8
# row 22
$ws.Cells.Item(22, 8).Value = 11250
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 14).ClearContents()
# row 31
$ws.Cells.Item(31, 8).Value = 30014
$ws.Cells.Item(31, 9).Value = 0
$ws.Cells.Item(31, 11).Value = 0
$ws.Cells.Item(31, 13).ClearContents()
# row 132
$ws.Cells.Item(132, 8).Value = 1279.7273
$ws.Cells.Item(132, 9).Value = 1074.4
$ws.Cells.Item(132, 11).Value = 3223.2
$ws.Cells.Item(132, 13).Value = -693.2000000000003
